# Scrum Activities.xlsx - "Add files via upload"
#
# The diff shows a new Daily Meeting row inserted above the existing
# "25/10/2018" row (old row 8, which becomes row 9), pushing every
# subsequent row down by one (old row 18 -> new row 19). The new row
# carries a date (as plain text, matching the surrounding "text date"
# rows already in the sheet) and a proof-file reference.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 8, shifting rows 8.. down by one.
$ws.Rows("8:8").Insert() | Out-Null

# Populate the newly inserted row with the new Scrum activity entry.
$ws.Range("A8").Value = "Daily Meeting"
$ws.Range("B8").Value = "15/11/2018"
$ws.Range("C8").Value = "Proof/SCRUMDAILY15-11-18.jpg"

# Match the author's resulting selection.
$ws.Range("C9").Select() | Out-Null
